$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timesheet")

# Row 96: new timesheet entry "Screencast maken"
$ws.Range("A96").Value = 43475
$ws.Range("B96").Value = "Screencast maken"
$ws.Range("C96").Value = 2

# Row 97: new timesheet entry "Presentatie maken"
$ws.Range("A97").Value = 43477
$ws.Range("B97").Value = "Presentatie maken"
$ws.Range("C97").Value = 1.5

# Update the selected cell to reflect where the user left off editing
$ws.Range("C98").Select()
